# Adds a new forecast-date column (R, "2020-04-19") and a new forecast
# row (30, "2020-05-03") to both the "cases" and "deaths" sheets, backfills
# the newly-observed "Observed" (column B) value for row 16, and corrects
# the sheet2 (deaths) B15 observed value.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cases")
$ws2 = $wb.Worksheets.Item("deaths")

$xlPasteValues = -4163

# --- column R header (new forecast date column), both sheets --------------
# R1 must hold the same text as A16 ("2020-04-19"), reusing the existing
# shared string rather than letting Excel reinterpret the literal text as
# a date serial. Copy + PasteSpecial(values) preserves the source cell's
# stored type instead of re-parsing text.
$ws1.Cells.Item(16, 1).Copy() | Out-Null
$ws1.Cells.Item(1, 18).PasteSpecial($xlPasteValues) | Out-Null
$ws2.Cells.Item(16, 1).Copy() | Out-Null
$ws2.Cells.Item(1, 18).PasteSpecial($xlPasteValues) | Out-Null

# --- materialize column R cells for rows 2-16 as empty (touched) cells ----
$ws1.Range("R2:R16").Style = "Normal"
$ws2.Range("R2:R16").Style = "Normal"

# --- column R forecast values for rows 17-29 -------------------------------
$ws1.Cells.Item(17, 18).Value2 = 40773
$ws1.Cells.Item(18, 18).Value2 = 43293
$ws1.Cells.Item(19, 18).Value2 = 46094
$ws1.Cells.Item(20, 18).Value2 = 50164
$ws1.Cells.Item(21, 18).Value2 = 52775
$ws1.Cells.Item(22, 18).Value2 = 55227
$ws1.Cells.Item(23, 18).Value2 = 57814
$ws1.Cells.Item(24, 18).Value2 = 61654
$ws1.Cells.Item(25, 18).Value2 = 64600
$ws1.Cells.Item(26, 18).Value2 = 68194
$ws1.Cells.Item(27, 18).Value2 = 73274
$ws1.Cells.Item(28, 18).Value2 = 76953
$ws1.Cells.Item(29, 18).Value2 = 82068

$ws2.Cells.Item(17, 18).Value2 = 2639
$ws2.Cells.Item(18, 18).Value2 = 2822
$ws2.Cells.Item(19, 18).Value2 = 3026
$ws2.Cells.Item(20, 18).Value2 = 3326
$ws2.Cells.Item(21, 18).Value2 = 3520
$ws2.Cells.Item(22, 18).Value2 = 3703
$ws2.Cells.Item(23, 18).Value2 = 3898
$ws2.Cells.Item(24, 18).Value2 = 4188
$ws2.Cells.Item(25, 18).Value2 = 4412
$ws2.Cells.Item(26, 18).Value2 = 4687
$ws2.Cells.Item(27, 18).Value2 = 5079
$ws2.Cells.Item(28, 18).Value2 = 5365
$ws2.Cells.Item(29, 18).Value2 = 5765

# --- backfill the "Observed" column (B) for row 16 -------------------------
$ws1.Cells.Item(16, 2).Value2 = 38654
$ws2.Cells.Item(16, 2).Value2 = 2462

# --- correct the row 15 "Observed" value on the deaths sheet ---------------
$ws2.Cells.Item(15, 2).Value2 = 2352

# --- new row 30 ("2020-05-03") ---------------------------------------------
# A30 needs the brand-new text "2020-05-03", which isn't anywhere else in
# the workbook yet. Build it via a formula in a scratch cell (so Excel
# treats the result as a plain string rather than auto-detecting a date),
# then copy the computed value across so no date number-format ever gets
# attached to a real cell.
$scratch = $ws1.Cells.Item(1, 26)
$scratch.Formula = '="2020-05-03"'
$scratch.Copy() | Out-Null
$ws1.Cells.Item(30, 1).PasteSpecial($xlPasteValues) | Out-Null
$ws2.Cells.Item(30, 1).PasteSpecial($xlPasteValues) | Out-Null

# Clear the scratch cell used to mint the new shared string.
$scratch.Clear() | Out-Null

# B30:Q30 are touched-but-empty on both sheets.
$ws1.Range("B30:Q30").Style = "Normal"
$ws2.Range("B30:Q30").Style = "Normal"

# R30 carries the new forecast value.
$ws1.Cells.Item(30, 18).Value2 = 84693
$ws2.Cells.Item(30, 18).Value2 = 5971

Write-Output "edit complete"
